$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlinks (on B2:B5) before rewriting the grid ---
$ws.Hyperlinks.Delete()

# --- Header row ---
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Firstname"
$ws.Range("C1").Value = "Lastname"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Address"
$ws.Range("F1").Value = "PhoneNumber"
$ws.Range("G1").Value = "Age"
$ws.Range("H1").Value = "Gender"

# --- Data rows ---
$ws.Range("A2").Value = "test1@gmail.com"
$ws.Range("B2").Value = "test1"
$ws.Range("C2").Value = "last1"
$ws.Range("D2").Value = "test123"
$ws.Range("E2").Value = "Hochiminh"
$ws.Range("F2").Value = 12345678
$ws.Range("G2").Value = 18
$ws.Range("H2").Value = "MALE"

$ws.Range("A3").Value = "test2@gmail.com"
$ws.Range("B3").Value = "test2"
$ws.Range("C3").Value = "last2"
$ws.Range("D3").Value = "test123"
$ws.Range("E3").Value = "Hochiminh"
$ws.Range("F3").Value = 12345678
$ws.Range("G3").Value = 18
$ws.Range("H3").Value = "FEMALE"

$ws.Range("A4").Value = "test3@gmail.com"
$ws.Range("B4").Value = "test3"
$ws.Range("C4").Value = "last3"
$ws.Range("D4").Value = "test123"
$ws.Range("E4").Value = "Hochiminh"
$ws.Range("F4").Value = 12345678
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = "OTHER"

$ws.Range("A5").Value = "test4@gmail.com"
$ws.Range("B5").Value = "test4"
$ws.Range("C5").Value = "last4"
$ws.Range("D5").Value = "test123"
$ws.Range("E5").Value = "Hochiminh"
$ws.Range("F5").Value = 12345678
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = "MALE"

# --- B2:H5 no longer carry the old "Hyperlink" cell style ---
$ws.Range("B2:B5").Style = "Normal"

# --- Hyperlinks on the email column (A2:A5); add first, style after so the
#     cells land on the same "Hyperlink" cellXf the sheet already defines ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:test2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:test3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:test4@gmail.com")
$ws.Range("A2:A5").Style = "Hyperlink"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.5
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- Selection matches the saved state ---
$ws.Range("H8").Select()

Write-Host "done"
